# "Added No Text CSV" -- fills in the missing Description column (D) on the
# "Model" sheet's feature-documentation table, re-labels the isVerified /
# isReply / isRetweeted rows (dropping the stale "Remove - None are
# retweets" note), and adds two new rows documenting retweet_count /
# favorite_count. Also touches a couple of cosmetic view/selection bits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# --- Make room for two new rows (retweet_count, favorite_count) right
# after the "user_since" row (row 13). Row 14 was already an empty spacer
# row, so only a single physical row needs inserting to end up with two
# free rows (14 & 15) before the "bigrams..." row, which will shift from
# row 15 down to row 16.
$ws.Rows.Item(15).Insert()

# --- Row 2: influence_interval ---------------------------------------
$ws.Cells.Item(2, 4).Value = "each interval = 2 hrs. add to created_hr to find time of day"

# --- Row 3: created_hr --------------------------------------------------
$ws.Cells.Item(3, 4).Value = "Midnight = 0"

# --- Row 4: dotw ---------------------------------------------------------
$ws.Cells.Item(4, 4).Value = "Sunday = 1"

# --- Row 5: num_tags -------------------------------------------------
$ws.Cells.Item(5, 4).Value = "number hashtags in tweet"

# --- Row 6: num_media ------------------------------------------------
$ws.Cells.Item(6, 4).Value = "number media in tweet"

# --- Row 7: num_symbols ------------------------------------------------
$ws.Cells.Item(7, 4).Value = "number symbols in tweet"

# --- Row 8: num_urls (description cell carries the black-font style) ---
$ws.Cells.Item(8, 4).Value = "number urls in tweet"
$ws.Cells.Item(8, 4).Font.Color = 0

# --- Row 9: num_user_mentions (same black-font style) -------------------
$ws.Cells.Item(9, 4).Value = "number users mentioned in tweet"
$ws.Cells.Item(9, 4).Font.Color = 0

# --- Row 10: was "isRetweeted" (red-highlighted, stale note) -> becomes
# "isVerified" with its real description. Clear the red highlight.
$ws.Cells.Item(10, 1).ClearFormats()
$ws.Cells.Item(10, 1).Value = "isVerified "
$ws.Cells.Item(10, 4).Value = "1 = the user is verified"

# --- Row 11: was "isVerified" -> becomes "isReply" ----------------------
$ws.Cells.Item(11, 1).Value = "isReply "
$ws.Cells.Item(11, 4).Value = "1 = the tweet is a reply to another tweet"

# --- Row 12: was "isReply" -> becomes "isRetweeted" (no trailing space) -
$ws.Cells.Item(12, 1).Value = "isRetweeted"
$ws.Cells.Item(12, 4).Value = "1 = the tweet has been retweeted "

# --- Row 13: user_since ---------------------------------------------
$ws.Cells.Item(13, 4).Value = "Number of years the user has been active"

# --- Row 14 (new): retweet_count -----------------------------------
$ws.Cells.Item(14, 1).Value = "retweet_count"
$ws.Cells.Item(14, 2).Value = "count"
$ws.Cells.Item(14, 3).Value = "0+"
$ws.Cells.Item(14, 4).Value = "number retweets"

# --- Row 15 (new): favorite_count ------------------------------------
$ws.Cells.Item(15, 1).Value = "favorite_count"
$ws.Cells.Item(15, 2).Value = "count"
$ws.Cells.Item(15, 3).Value = "0+"
$ws.Cells.Item(15, 4).Value = "number favorites"

# --- Column D is wider now that it holds real descriptions --------------
$ws.Columns.Item(4).ColumnWidth = 58

# --- Selection cosmetics, matching where the author left the cursor -----
$ws.Range("D22").Select()

$ws1 = $wb.Worksheets.Item("Determine Unluence")
$ws1.Range("D34").Select()

$ws.Activate()
